$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new metadata column (metadatafield5) was added to the sheet, along
# with uppercase-test values for the two existing sample rows.

# New header cell for the 5th metadata field - same style as the rest
# of the header row.
$ws.Range("F1").Value = "metadatafield5"
$ws.Range("F1").Style = $ws.Range("D1").Style

# The header row now shares one common (plain) style across all cells -
# realign E1 with the rest of the header row instead of its previous,
# distinct style.
$ws.Range("E1").Style = $ws.Range("D1").Style

# New metadata values for the existing sample rows.
$ws.Range("E2").Value = "A Test"
$ws.Range("E3").Value = "Another Test"

# The new text values get an explicit text number format.
$ws.Range("E2:E3").NumberFormat = "@"

# Matches the cursor position recorded after the edit.
$ws.Range("E2").Select()
